# This script appends newly-logged sensor readings (2026-01-28, ~16:22-16:23)
# to the PIR, Humidity, and Temperature sheets of the SeniorConnect master log.
$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append new sensor log rows ----
$ws = $wb.Worksheets.Item("PIR")
$newRows = @(
    @("2026-01-28", "16:22:40", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:22:42", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:22:47", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:22:52", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:22:57", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:02", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:07", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:12", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:17", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:22", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:27", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:32", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:23:37", "16:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    # Column A holds a literal date-like string ("YYYY-MM-DD"); force text
    # formatting first so Excel does not auto-convert it into a date serial.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

# ---- Humidity sheet: append new sensor log rows ----
$ws = $wb.Worksheets.Item("Humidity")
$newRows = @(
    @("2026-01-28", "16:22:39", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-28", "16:22:40", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:22:41", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-28", "16:22:45", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:22:49", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:22:53", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-28", "16:22:57", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:23:01", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-28", "16:23:09", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:23:13", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:23:17", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:23:25", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:23:29", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:23:33", "16:00", "Bathroom", "87.3%", "Active"),
    @("2026-01-28", "16:23:37", "16:00", "Bathroom", "88.2%", "Active")
)
$startRow = 144
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    # Column A holds a literal date-like string ("YYYY-MM-DD"); force text
    # formatting first so Excel does not auto-convert it into a date serial.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    # Column E holds a literal percentage string (e.g. "87.3%"); force text
    # formatting first so Excel does not auto-convert it into a percentage number.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

# ---- Temperature sheet: append new sensor log rows ----
$ws = $wb.Worksheets.Item("Temperature")
$newRows = @(
    @("2026-01-28", "16:22:39", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:40", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:41", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:45", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:49", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:53", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:22:58", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:02", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:10", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:13", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:18", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:26", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:30", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:34", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:23:38", "16:00", "Bathroom", "22.8C", "Active")
)
$startRow = 144
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    # Column A holds a literal date-like string ("YYYY-MM-DD"); force text
    # formatting first so Excel does not auto-convert it into a date serial.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
